$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 45 new data rows (sheet rows 102-146) continuing the existing
# regcntr_id / machine_id / device_id test-data pattern:
#   A cycles 10002..10010 (9-row period), B cycles 10021..10029 (9-row
#   period, same phase as A), C increments by 1 starting at 3000121.
# D/E/F/G/H stay constant for every data row in this sheet
# ("eng", TRUE, "superadmin()", "now()", "now()").
for ($i = 0; $i -lt 45; $i++) {
    $row = 102 + $i
    $a = 10002 + ($i % 9)
    $b = 10021 + ($i % 9)
    $c = 3000121 + $i

    $ws.Cells.Item($row, 1).Value = $a
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $True
    $ws.Cells.Item($row, 6).Value = "superadmin()"
    $ws.Cells.Item($row, 7).Value = "now()"
    $ws.Cells.Item($row, 8).Value = "now()"
}

# The saved file's cursor/selection ends up parked just past the new
# data (row 147, full-column selection) - mirror that.
$ws.Range("A147:XFD1048576").Select()

# The sheet was set up for printing (portrait) when it was resaved.
$ws.PageSetup.Orientation = 1
